$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 299
$ws.Range("F6").Value = 430
$ws.Range("F7").Value = 371
$ws.Range("F8").Value = 1956
$ws.Range("F12").Value = 1586
$ws.Range("F13").Value = 1586
$ws.Range("F15").Value = 47
$ws.Range("F24").Value = 6926
$ws.Range("F25").Value = 7494
$ws.Range("F29").Value = 50
$ws.Range("F31").Value = 235
$ws.Range("F39").Value = 274
$ws.Range("F40").Value = 670
$ws.Range("F43").Value = 304
$ws.Range("F45").Value = 179
$ws.Range("F46").Value = 77
$ws.Range("F47").Value = 115
$ws.Range("F48").Value = 130

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 21

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 107

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 107
$ws.Range("F9").Value = 299
$ws.Range("F11").Value = 430
$ws.Range("F12").Value = 371
$ws.Range("F13").Value = 1956
$ws.Range("F16").Value = 1586
$ws.Range("F17").Value = 1586
$ws.Range("F18").Value = 47
$ws.Range("F24").Value = 6926
$ws.Range("F25").Value = 7494
$ws.Range("F27").Value = 235
$ws.Range("F34").Value = 274
$ws.Range("F35").Value = 21
$ws.Range("F37").Value = 670
$ws.Range("F43").Value = 304
$ws.Range("F45").Value = 179
$ws.Range("F46").Value = 77
$ws.Range("F47").Value = 115
